# Saldo_guide.xlsx update: refresh reference date (Dt. Referencia) for all
# data rows from 2024-04-22 (45404) to 2024-04-23 (45405), and correct a
# handful of Saldo Previsto / Vl. Projetado values that were recalculated
# for the new reference date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Dt. Referencia" column (G) for every data row (2-310) by one day.
$ws.Range("G2:G310").Value = 45405

# Row 55 - FELIPE ROSSI RAMOS: Saldo Previsto now matches Vl. Total, Vl. Projetado zeroed.
$ws.Range("D55").Value = 16205.42
$ws.Range("E55").Value = 0

# Row 76: Saldo Previsto and Vl. Total both zeroed out.
$ws.Range("D76").Value = 0
$ws.Range("H76").Value = 0

# Row 129: Saldo Previsto / Vl. Total updated to new balance.
$ws.Range("D129").Value = 88.88
$ws.Range("H129").Value = 88.88

# Row 135: Saldo Previsto / Vl. Total updated to new balance.
$ws.Range("D135").Value = 12234.67
$ws.Range("H135").Value = 12234.67

# Row 290: Saldo Previsto now matches Vl. Total, Vl. Projetado zeroed.
$ws.Range("D290").Value = 46785.7
$ws.Range("E290").Value = 0
